# Auto-generated edit script: update the worksheet date and all 100
# addition/subtraction answer cells to the new values from the commit.
$d = $word.ActiveDocument

$replacements = @(
    @("2025-06-14 Saturday", "2025-06-15 Sunday"),
    @("86-26=60", "74+6=80"),
    @("56-18=38", "5+70=75"),
    @("24-3=21", "78-32=46"),
    @("0+16=16", "28+54=82"),
    @("79-31=48", "22+12=34"),
    @("71+25=96", "15-0=15"),
    @("45-6=39", "30+36=66"),
    @("19+30=49", "91+0=91"),
    @("58+21=79", "76-0=76"),
    @("6+4=10", "91-80=11"),
    @("68-50=18", "44-41=3"),
    @("24-4=20", "12+67=79"),
    @("19+7=26", "28+35=63"),
    @("78-42=36", "67-52=15"),
    @("72-1=71", "69-18=51"),
    @("86-68=18", "17+19=36"),
    @("72-0=72", "46+34=80"),
    @("74-73=1", "73+9=82"),
    @("23+69=92", "46+52=98"),
    @("47+6=53", "48-35=13"),
    @("82+10=92", "56-13=43"),
    @("61+28=89", "41-14=27"),
    @("35-28=7", "9+61=70"),
    @("45-8=37", "3+21=24"),
    @("42-40=2", "44+45=89"),
    @("65+21=86", "28+70=98"),
    @("26+68=94", "51-32=19"),
    @("16+14=30", "31+17=48"),
    @("37+23=60", "68-62=6"),
    @("78-1=77", "57-50=7"),
    @("75-8=67", "82+5=87"),
    @("95-32=63", "83-32=51"),
    @("79-38=41", "0+68=68"),
    @("94-89=5", "81+0=81"),
    @("96-43=53", "63+9=72"),
    @("0+31=31", "95-24=71"),
    @("63+17=80", "38-17=21"),
    @("94-73=21", "44-34=10"),
    @("14+37=51", "27+61=88"),
    @("13+48=61", "75-42=33"),
    @("45+38=83", "1+15=16"),
    @("61-25=36", "92-58=34"),
    @("7+23=30", "15+42=57"),
    @("31-27=4", "25+5=30"),
    @("54+5=59", "88-54=34"),
    @("2+2=4", "3+72=75"),
    @("1+66=67", "45+20=65"),
    @("44-19=25", "20+61=81"),
    @("68-53=15", "97-68=29"),
    @("23+24=47", "12-0=12"),
    @("64+29=93", "95-59=36"),
    @("9+64=73", "67+15=82"),
    @("11+14=25", "74+19=93"),
    @("20+71=91", "59-26=33"),
    @("85-54=31", "49+18=67"),
    @("76-51=25", "60+1=61"),
    @("25-13=12", "39-0=39"),
    @("69-29=40", "72-56=16"),
    @("72+3=75", "12+21=33"),
    @("4+61=65", "99-23=76"),
    @("42+48=90", "0+55=55"),
    @("18+68=86", "3+63=66"),
    @("71-39=32", "45-3=42"),
    @("47+26=73", "10+43=53"),
    @("69-28=41", "85-44=41"),
    @("29+52=81", "27+30=57"),
    @("51+9=60", "5-1=4"),
    @("77-11=66", "90+1=91"),
    @("48+7=55", "38-3=35"),
    @("98-8=90", "85-35=50"),
    @("34-5=29", "28+33=61"),
    @("26+7=33", "26+1=27"),
    @("99-54=45", "31+26=57"),
    @("10+71=81", "90-0=90"),
    @("35+56=91", "63-0=63"),
    @("34+41=75", "13+1=14"),
    @("38+31=69", "23+50=73"),
    @("19+25=44", "10-5=5"),
    @("32-14=18", "8+31=39"),
    @("58-29=29", "47+24=71"),
    @("65+9=74", "95-73=22"),
    @("95-77=18", "7+31=38"),
    @("47+37=84", "3-0=3"),
    @("37+43=80", "33+14=47"),
    @("68-58=10", "84-16=68"),
    @("14+13=27", "81-13=68"),
    @("22-0=22", "10+27=37"),
    @("35+2=37", "87-50=37"),
    @("24+0=24", "57-2=55"),
    @("21+9=30", "47+23=70"),
    @("85+6=91", "74-8=66"),
    @("68-11=57", "88-64=24"),
    @("76-27=49", "63+25=88"),
    @("52-18=34", "75-65=10"),
    @("41+23=64", "42+57=99"),
    @("54+25=79", "85-68=17"),
    @("27+7=34", "13+55=68"),
    @("48+12=60", "95-5=90"),
    @("5+14=19", "66-55=11"),
    @("24-15=9", "40-37=3"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

